$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018")

# KW52 (row 57) vacation count goes from 1 to 0.
# Dependent formulas (B2 = B5/B3, B4 = B3-B5, B5 = SUM(B6:B58)) will
# recalculate automatically because of this change.
$ws.Range("B57").Value = 0
